$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: DAMSLTag (I12) / DialogAct (J12)
$ws.Range("I12").Value = "aa"
$ws.Range("J12").Value = "Agree/Accept"

# Row 19: DAMSLTag (I19) / DialogAct (J19)
$ws.Range("I19").Value = "ba"
$ws.Range("J19").Value = "Appreciation"

# Row 25: DAMSLTag (I25) / DialogAct (J25)
$ws.Range("I25").Value = "sv"
$ws.Range("J25").Value = "Statement-opinion"
